$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 187.64706
$ws.Range("I2").Value = 187.64706
$ws.Range("K2").Value = 187.64706
$ws.Range("M2").Value = -74.64706000000001
$ws.Range("H5").Value = 69.8
$ws.Range("I5").Value = 69.8
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 69.8
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 45.2
$ws.Range("N5").ClearContents()
$ws.Range("H32").Value = 1500
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1500
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1500
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -2152
$ws.Range("H40").Value = 123557.64
$ws.Range("I40").Value = 751235.25
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 751235.25
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -751060.25
$ws.Range("N40").Value = -4350
$ws.Range("H44").Value = 45997.5
$ws.Range("J44").Value = 45997.5
$ws.Range("L44").Value = 45997.5
$ws.Range("N44").Value = -46921.5
$ws.Range("H53").Value = 1473.8422
$ws.Range("I53").Value = 1019.6
$ws.Range("J53").Value = 1978.5555
$ws.Range("K53").Value = 1019.6
$ws.Range("L53").Value = 1978.5555
$ws.Range("M53").Value = -382.6
$ws.Range("N53").Value = -3252.5555
$ws.Range("H62").Value = 8583.538
$ws.Range("I62").Value = 7433.5557
$ws.Range("K62").Value = 7433.5557
$ws.Range("M62").Value = -6809.5557
$ws.Range("H65").Value = 8583.538
$ws.Range("I65").Value = 7433.5557
$ws.Range("K65").Value = 37167.7785
$ws.Range("M65").Value = -34047.7785
$ws.Range("H94").Value = 1852.5714
$ws.Range("I94").Value = 1852.5714
$ws.Range("K94").Value = 1852.5714
$ws.Range("M94").Value = -1401.5714
$ws.Range("H98").Value = 1690.2632
$ws.Range("I98").Value = 1354.125
$ws.Range("K98").Value = 1354.125
$ws.Range("M98").Value = 143.875
$ws.Range("H106").Value = 4169322.5
$ws.Range("I106").Value = 4447204
$ws.Range("J106").Value = 1099
$ws.Range("K106").Value = 4447204
$ws.Range("L106").Value = 1099
$ws.Range("M106").Value = -4446573
$ws.Range("N106").Value = -2361
$ws.Range("H107").Value = 1756.2954
$ws.Range("I107").Value = 1796.85
$ws.Range("K107").Value = 1796.85
$ws.Range("M107").Value = 123.1500000000001
$ws.Range("H122").Value = 1690.2632
$ws.Range("I122").Value = 1354.125
$ws.Range("K122").Value = 4062.375
$ws.Range("M122").Value = -1612.375
$ws.Range("H132").Value = 3096.3606
$ws.Range("I132").Value = 3114.6333
$ws.Range("K132").Value = 9343.8999
$ws.Range("M132").Value = -6813.8999
$ws.Range("H135").Value = 1096.0714
$ws.Range("I135").Value = 1141.9615
$ws.Range("K135").Value = 10277.6535
$ws.Range("M135").Value = -7742.653499999999
$ws.Range("H137").Value = 43485.45
$ws.Range("I137").Value = 50938.47
$ws.Range("J137").Value = 1251.6666
$ws.Range("K137").Value = 152815.41
$ws.Range("L137").Value = 3754.9998
$ws.Range("M137").Value = -150265.41
$ws.Range("N137").Value = -8854.9998
$ws.Range("H138").Value = 3158.03
$ws.Range("I138").Value = 1244.0769
$ws.Range("J138").Value = 3830.5
$ws.Range("K138").Value = 3732.2307
$ws.Range("L138").Value = 11491.5
$ws.Range("M138").Value = 1407.7693
$ws.Range("N138").Value = -21771.5
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 300
$ws.Range("K4").Value = 300
$ws.Range("M4").Value = -184
$ws.Range("H32").Value = 12511336
$ws.Range("I32").Value = 13031881
$ws.Range("K32").Value = 13031881
$ws.Range("M32").Value = -13031594
$ws.Range("H45").Value = 3549.72
$ws.Range("I45").Value = 3151.111
$ws.Range("K45").Value = 3151.111
$ws.Range("M45").Value = -2774.111
$ws.Range("H61").Value = 6549.846
$ws.Range("I61").Value = 8056.125
$ws.Range("J61").Value = 4139.8
$ws.Range("K61").Value = 8056.125
$ws.Range("L61").Value = 4139.8
$ws.Range("M61").Value = -7844.125
$ws.Range("N61").Value = -4563.8
$ws.Range("H63").Value = 4685.7144
$ws.Range("J63").Value = 5650
$ws.Range("L63").Value = 5650
$ws.Range("N63").Value = -7022
$ws.Range("H66").Value = 4685.7144
$ws.Range("J66").Value = 5650
$ws.Range("L66").Value = 28250
$ws.Range("N66").Value = -35114
$ws.Range("H74").Value = 3580.8
$ws.Range("I74").Value = 2906.3157
$ws.Range("J74").Value = 5716.6665
$ws.Range("K74").Value = 2906.3157
$ws.Range("L74").Value = 5716.6665
$ws.Range("M74").Value = -2032.3157
$ws.Range("N74").Value = -7464.6665
$ws.Range("H77").Value = 3580.8
$ws.Range("I77").Value = 2906.3157
$ws.Range("J77").Value = 5716.6665
$ws.Range("K77").Value = 14531.5785
$ws.Range("L77").Value = 28583.3325
$ws.Range("M77").Value = -10163.5785
$ws.Range("N77").Value = -37319.3325
$ws.Range("H97").Value = 531.2
$ws.Range("I97").Value = 373.55554
$ws.Range("K97").Value = 373.55554
$ws.Range("M97").Value = 122.44446
$ws.Range("H102").Value = 1785.3846
$ws.Range("I102").Value = 901.2222
$ws.Range("J102").Value = 3774.75
$ws.Range("K102").Value = 901.2222
$ws.Range("L102").Value = 3774.75
$ws.Range("M102").Value = 720.7778
$ws.Range("N102").Value = -7018.75
$ws.Range("H108").Value = 79191.60000000001
$ws.Range("J108").Value = 79191.60000000001
$ws.Range("L108").Value = 79191.60000000001
$ws.Range("N108").Value = -86871.60000000001
$ws.Range("H110").Value = 2493.923
$ws.Range("I110").Value = 2142.2
$ws.Range("K110").Value = 2142.2
$ws.Range("M110").Value = -97.19999999999982
$ws.Range("H132").Value = 3790.9678
$ws.Range("I132").Value = 3568.7
$ws.Range("K132").Value = 10706.1
$ws.Range("M132").Value = -8176.099999999999
$ws.Range("H136").Value = 6549.846
$ws.Range("I136").Value = 8056.125
$ws.Range("J136").Value = 4139.8
$ws.Range("K136").Value = 24168.375
$ws.Range("L136").Value = 12419.4
$ws.Range("M136").Value = -21618.375
$ws.Range("N136").Value = -17519.4
$ws.Range("H138").Value = 99999.336
$ws.Range("I138").Value = 100000
$ws.Range("J138").Value = 99999
$ws.Range("K138").Value = 100000
$ws.Range("L138").Value = 99999
$ws.Range("M138").Value = -94860
$ws.Range("N138").Value = -110279
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 18809.104
$ws.Range("I20").Value = 25006.715
$ws.Range("J20").Value = 2540.375
$ws.Range("K20").Value = 25006.715
$ws.Range("L20").Value = 2540.375
$ws.Range("M20").Value = -24759.715
$ws.Range("N20").Value = -3034.375
$ws.Range("H86").Value = 2662.8333
$ws.Range("I86").Value = 1998.3334
$ws.Range("J86").Value = 3327.3333
$ws.Range("K86").Value = 1998.3334
$ws.Range("L86").Value = 3327.3333
$ws.Range("M86").Value = -875.3334
$ws.Range("N86").Value = -5573.3333
$ws.Range("H89").Value = 2662.8333
$ws.Range("I89").Value = 1998.3334
$ws.Range("J89").Value = 3327.3333
$ws.Range("K89").Value = 9991.666999999999
$ws.Range("L89").Value = 16636.6665
$ws.Range("M89").Value = -4375.666999999999
$ws.Range("N89").Value = -27868.6665
$ws.Range("H99").Value = 2586.037
$ws.Range("I99").Value = 1319
$ws.Range("J99").Value = 3950.5386
$ws.Range("K99").Value = 1319
$ws.Range("L99").Value = 3950.5386
$ws.Range("M99").Value = 179
$ws.Range("N99").Value = -6946.5386
$ws.Range("H105").Value = 2474
$ws.Range("I105").Value = 2161.15
$ws.Range("J105").Value = 3725.4
$ws.Range("K105").Value = 2161.15
$ws.Range("L105").Value = 3725.4
$ws.Range("M105").Value = -414.1500000000001
$ws.Range("N105").Value = -7219.4
$ws.Range("H107").Value = 1783.2858
$ws.Range("I107").Value = 1462.8182
$ws.Range("J107").Value = 2958.3333
$ws.Range("K107").Value = 1462.8182
$ws.Range("L107").Value = 2958.3333
$ws.Range("M107").Value = 457.1818000000001
$ws.Range("N107").Value = -6798.3333
$ws.Range("H134").Value = 3249350.5
$ws.Range("I134").Value = 4466218
$ws.Range("J134").Value = 4371
$ws.Range("K134").Value = 13398654
$ws.Range("L134").Value = 13113
$ws.Range("M134").Value = -13396119
$ws.Range("N134").Value = -18183
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 836.125
$ws.Range("I22").Value = 824
$ws.Range("K22").Value = 824
$ws.Range("M22").Value = -474
$ws.Range("I31").Value = 2071.95
$ws.Range("J31").Value = 4465.7334
$ws.Range("K31").Value = 2071.95
$ws.Range("L31").Value = 4465.7334
$ws.Range("M31").Value = -1776.95
$ws.Range("N31").Value = -5055.7334
$ws.Range("I34").Value = 2071.95
$ws.Range("J34").Value = 4465.7334
$ws.Range("K34").Value = 2071.95
$ws.Range("L34").Value = 4465.7334
$ws.Range("M34").Value = -1869.95
$ws.Range("N34").Value = -4869.7334
$ws.Range("H53").Value = 26499.5
$ws.Range("J53").Value = 26999
$ws.Range("L53").Value = 26999
$ws.Range("N53").Value = -28213
$ws.Range("H58").Value = 3191.261
$ws.Range("I58").Value = 3190.4
$ws.Range("K58").Value = 3190.4
$ws.Range("M58").Value = -2987.4
$ws.Range("H60").Value = 11666.667
$ws.Range("I60").Value = 11666.667
$ws.Range("K60").Value = 11666.667
$ws.Range("M60").Value = -11155.667
$ws.Range("H69").Value = 69095.625
$ws.Range("I69").Value = 55666.332
$ws.Range("J69").Value = 77153.2
$ws.Range("K69").Value = 55666.332
$ws.Range("L69").Value = 77153.2
$ws.Range("M69").Value = -54917.332
$ws.Range("N69").Value = -78651.2
$ws.Range("H72").Value = 69095.625
$ws.Range("I72").Value = 55666.332
$ws.Range("J72").Value = 77153.2
$ws.Range("K72").Value = 166998.996
$ws.Range("L72").Value = 231459.6
$ws.Range("M72").Value = -163254.996
$ws.Range("N72").Value = -238947.6
$ws.Range("H107").Value = 86132.914
$ws.Range("I107").Value = 253574.25
$ws.Range("K107").Value = 253574.25
$ws.Range("M107").Value = -251654.25
$ws.Range("H132").Value = 56895.297
$ws.Range("I132").Value = 33127.637
$ws.Range("K132").Value = 99382.91100000001
$ws.Range("M132").Value = -96852.91100000001
$ws.Range("H134").Value = 2367.9
$ws.Range("I134").Value = 1972.375
$ws.Range("J134").Value = 3950
$ws.Range("K134").Value = 5917.125
$ws.Range("L134").Value = 11850
$ws.Range("M134").Value = -3382.125
$ws.Range("N134").Value = -16920
$ws.Range("H136").Value = 3191.261
$ws.Range("I136").Value = 3190.4
$ws.Range("K136").Value = 9571.200000000001
$ws.Range("M136").Value = -7021.200000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1417.3334
$ws.Range("I5").Value = 1417.3334
$ws.Range("K5").Value = 4252.0002
$ws.Range("M5").Value = -4140.0002
$ws.Range("H7").Value = 840.03125
$ws.Range("I7").Value = 573.4167
$ws.Range("K7").Value = 1720.2501
$ws.Range("M7").Value = -1608.2501
$ws.Range("H8").Value = 156.33333
$ws.Range("I8").Value = 156.33333
$ws.Range("K8").Value = 468.99999
$ws.Range("M8").Value = -329.99999
$ws.Range("H26").Value = 160
$ws.Range("I26").Value = 20
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 60
$ws.Range("L26").Value = 900
$ws.Range("M26").Value = 228
$ws.Range("N26").Value = -1476
$ws.Range("H47").Value = 742.25
$ws.Range("I47").Value = 789.6667
$ws.Range("J47").Value = 600
$ws.Range("K47").Value = 2369.0001
$ws.Range("L47").Value = 1800
$ws.Range("M47").Value = -1938.0001
$ws.Range("N47").Value = -2662
$ws.Range("H80").Value = 4383.846
$ws.Range("I80").Value = 3333
$ws.Range("J80").Value = 4699.1
$ws.Range("K80").Value = 9999
$ws.Range("L80").Value = 14097.3
$ws.Range("M80").Value = -9063
$ws.Range("N80").Value = -15969.3
$ws.Range("H83").Value = 4383.846
$ws.Range("I83").Value = 3333
$ws.Range("J83").Value = 4699.1
$ws.Range("K83").Value = 29997
$ws.Range("L83").Value = 42291.9
$ws.Range("M83").Value = -25317
$ws.Range("N83").Value = -51651.9
$ws.Range("H92").Value = 1465.6666
$ws.Range("I92").Value = 1498
$ws.Range("J92").Value = 1449.5
$ws.Range("K92").Value = 4494
$ws.Range("L92").Value = 4348.5
$ws.Range("M92").Value = -3246
$ws.Range("N92").Value = -6844.5
$ws.Range("H132").Value = 1022.8214
$ws.Range("I132").Value = 1046.1305
$ws.Range("K132").Value = 9415.174499999999
$ws.Range("M132").Value = -6885.174499999999
$ws.Range("H135").Value = 1417.3334
$ws.Range("I135").Value = 1417.3334
$ws.Range("K135").Value = 12756.0006
$ws.Range("M135").Value = -10221.0006
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2191.6
$ws.Range("I80").Value = 1739.75
$ws.Range("J80").Value = 3999
$ws.Range("K80").Value = 1739.75
$ws.Range("L80").Value = 3999
$ws.Range("M80").Value = -741.75
$ws.Range("N80").Value = -5995
$ws.Range("H83").Value = 2191.6
$ws.Range("I83").Value = 1739.75
$ws.Range("J83").Value = 3999
$ws.Range("K83").Value = 8698.75
$ws.Range("L83").Value = 19995
$ws.Range("M83").Value = -3706.75
$ws.Range("N83").Value = -29979
$ws.Range("H113").Value = 689.8
$ws.Range("I113").Value = 599.75
$ws.Range("J113").Value = 1050
$ws.Range("K113").Value = 599.75
$ws.Range("L113").Value = 1050
$ws.Range("M113").Value = 1570.25
$ws.Range("N113").Value = -5390
$ws.Range("H122").Value = 3666.3
$ws.Range("I122").Value = 4057
$ws.Range("K122").Value = 12171
$ws.Range("M122").Value = -9721
$ws.Range("H126").Value = 3766.1667
$ws.Range("I126").Value = 4033
$ws.Range("J126").Value = 3499.3333
$ws.Range("K126").Value = 12099
$ws.Range("L126").Value = 10497.9999
$ws.Range("M126").Value = -9629
$ws.Range("N126").Value = -15437.9999
$ws.Range("H132").Value = 3948
$ws.Range("I132").Value = 3853.2188
$ws.Range("J132").Value = 4554.6
$ws.Range("K132").Value = 11559.6564
$ws.Range("L132").Value = 13663.8
$ws.Range("M132").Value = -9029.6564
$ws.Range("N132").Value = -18723.8
$ws.Range("H135").Value = 300000
$ws.Range("J135").Value = 300000
$ws.Range("L135").Value = 300000
$ws.Range("N135").Value = -310140
$ws.Range("H136").Value = 69341.86
$ws.Range("J136").Value = 72732.164
$ws.Range("L136").Value = 218196.492
$ws.Range("N136").Value = -223296.492
$ws.Range("H137").Value = 40000
$ws.Range("I137").Value = 40000
$ws.Range("K137").Value = 40000
$ws.Range("M137").Value = -34900
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H140").Value = 291999.5
$ws.Range("J140").Value = 291999.5
$ws.Range("L140").Value = 291999.5
$ws.Range("N140").Value = -302359.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8452.299999999999
$ws.Range("I7").Value = 10712.375
$ws.Range("J7").Value = 6945.5835
$ws.Range("K7").Value = 10712.375
$ws.Range("L7").Value = 6945.5835
$ws.Range("M7").Value = -10600.375
$ws.Range("N7").Value = -7169.5835
$ws.Range("H40").Value = 25648178
$ws.Range("I40").Value = 41671916
$ws.Range("J40").Value = 10197.4
$ws.Range("K40").Value = 41671916
$ws.Range("L40").Value = 10197.4
$ws.Range("M40").Value = -41671780
$ws.Range("N40").Value = -10469.4
$ws.Range("H46").Value = 6656.6
$ws.Range("I46").Value = 8945.866
$ws.Range("J46").Value = 4367.3335
$ws.Range("K46").Value = 8945.866
$ws.Range("L46").Value = 4367.3335
$ws.Range("M46").Value = -8757.866
$ws.Range("N46").Value = -4743.3335
$ws.Range("H55").Value = 436.31818
$ws.Range("I55").Value = 235.90909
$ws.Range("J55").Value = 636.7273
$ws.Range("K55").Value = 235.90909
$ws.Range("L55").Value = 636.7273
$ws.Range("M55").Value = -62.90908999999999
$ws.Range("N55").Value = -982.7273
$ws.Range("H56").Value = 6081.4
$ws.Range("I56").Value = 6081.4
$ws.Range("K56").Value = 6081.4
$ws.Range("M56").Value = -5390.4
$ws.Range("H58").Value = 7949.5
$ws.Range("I58").Value = 7949.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 7949.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -7689.5
$ws.Range("N58").ClearContents()
$ws.Range("H61").Value = 2669.3333
$ws.Range("I61").Value = 2909.25
$ws.Range("J61").Value = 750
$ws.Range("K61").Value = 2909.25
$ws.Range("L61").Value = 750
$ws.Range("M61").Value = -2707.25
$ws.Range("N61").Value = -1154
$ws.Range("H68").Value = 3974.0264
$ws.Range("I68").Value = 3250.3333
$ws.Range("K68").Value = 3250.3333
$ws.Range("M68").Value = -2501.3333
$ws.Range("H71").Value = 3974.0264
$ws.Range("I71").Value = 3250.3333
$ws.Range("K71").Value = 16251.6665
$ws.Range("M71").Value = -12507.6665
$ws.Range("H82").Value = 2060.6428
$ws.Range("I82").Value = 2094.4443
$ws.Range("J82").Value = 1999.8
$ws.Range("K82").Value = 2094.4443
$ws.Range("L82").Value = 1999.8
$ws.Range("M82").Value = -1733.4443
$ws.Range("N82").Value = -2721.8
$ws.Range("H85").Value = 2060.6428
$ws.Range("I85").Value = 2094.4443
$ws.Range("J85").Value = 1999.8
$ws.Range("K85").Value = 2094.4443
$ws.Range("L85").Value = 1999.8
$ws.Range("M85").Value = -846.4443000000001
$ws.Range("N85").Value = -4495.8
$ws.Range("H113").Value = 2669.3333
$ws.Range("I113").Value = 2909.25
$ws.Range("J113").Value = 750
$ws.Range("K113").Value = 2909.25
$ws.Range("L113").Value = 750
$ws.Range("M113").Value = -739.25
$ws.Range("N113").Value = -5090
$ws.Range("H126").Value = 8452.299999999999
$ws.Range("I126").Value = 10712.375
$ws.Range("J126").Value = 6945.5835
$ws.Range("K126").Value = 32137.125
$ws.Range("L126").Value = 20836.7505
$ws.Range("M126").Value = -29667.125
$ws.Range("N126").Value = -25776.7505
$ws.Range("H132").Value = 90131.46000000001
$ws.Range("I132").Value = 104754.45
$ws.Range("K132").Value = 314263.35
$ws.Range("M132").Value = -311733.35
$ws.Range("H133").Value = 67471.25
$ws.Range("I133").Value = 90001
$ws.Range("J133").Value = 44941.5
$ws.Range("K133").Value = 90001
$ws.Range("L133").Value = 44941.5
$ws.Range("M133").Value = -87471
$ws.Range("N133").Value = -50001.5
$ws.Range("H136").Value = 2362.0557
$ws.Range("I136").Value = 1617
$ws.Range("J136").Value = 3293.375
$ws.Range("K136").Value = 4851
$ws.Range("L136").Value = 9880.125
$ws.Range("M136").Value = -2301
$ws.Range("N136").Value = -14980.125
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 11110.2
$ws.Range("I32").Value = 11110.2
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 11110.2
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -10793.2
$ws.Range("N32").ClearContents()
$ws.Range("H33").Value = 19000
$ws.Range("I33").Value = 19000
$ws.Range("K33").Value = 19000
$ws.Range("M33").Value = -18750
$ws.Range("H36").Value = 19000
$ws.Range("I36").Value = 19000
$ws.Range("K36").Value = 19000
$ws.Range("M36").Value = -18750
$ws.Range("H51").Value = 42365.668
$ws.Range("I51").Value = 24733.334
$ws.Range("J51").Value = 59998
$ws.Range("K51").Value = 24733.334
$ws.Range("L51").Value = 59998
$ws.Range("M51").Value = -24223.334
$ws.Range("N51").Value = -61018
$ws.Range("H54").Value = 56267.5
$ws.Range("I54").Value = 30070
$ws.Range("J54").Value = 65000
$ws.Range("K54").Value = 30070
$ws.Range("L54").Value = 65000
$ws.Range("M54").Value = -29550
$ws.Range("N54").Value = -66040
$ws.Range("H70").Value = 43821
$ws.Range("I70").Value = 73000
$ws.Range("J70").Value = 37985.2
$ws.Range("K70").Value = 73000
$ws.Range("L70").Value = 37985.2
$ws.Range("M70").Value = -72685
$ws.Range("N70").Value = -38615.2
$ws.Range("H73").Value = 43821
$ws.Range("I73").Value = 73000
$ws.Range("J73").Value = 37985.2
$ws.Range("K73").Value = 73000
$ws.Range("L73").Value = 37985.2
$ws.Range("M73").Value = -71908
$ws.Range("N73").Value = -40169.2
$ws.Range("H81").Value = 4339.4
$ws.Range("I81").Value = 3898.9167
$ws.Range("J81").Value = 5000.125
$ws.Range("K81").Value = 7797.8334
$ws.Range("L81").Value = 10000.25
$ws.Range("M81").Value = -6736.8334
$ws.Range("N81").Value = -12122.25
$ws.Range("H84").Value = 4339.4
$ws.Range("I84").Value = 3898.9167
$ws.Range("J84").Value = 5000.125
$ws.Range("K84").Value = 38989.167
$ws.Range("L84").Value = 50001.25
$ws.Range("M84").Value = -33685.167
$ws.Range("N84").Value = -60609.25
$ws.Range("H92").Value = 66244.5
$ws.Range("J92").Value = 66244.5
$ws.Range("L92").Value = 66244.5
$ws.Range("N92").Value = -71236.5
$ws.Range("H96").Value = 10049.25
$ws.Range("I96").Value = 7048.5
$ws.Range("K96").Value = 7048.5
$ws.Range("M96").Value = -5675.5
$ws.Range("H113").Value = 1912.5
$ws.Range("I113").Value = 1758.3334
$ws.Range("J113").Value = 2375
$ws.Range("K113").Value = 5275.0002
$ws.Range("L113").Value = 7125
$ws.Range("M113").Value = -3105.0002
$ws.Range("N113").Value = -11465
$ws.Range("H122").Value = 66669976
$ws.Range("J122").Value = 4917.5557
$ws.Range("L122").Value = 14752.6671
$ws.Range("N122").Value = -19652.6671
$ws.Range("H126").Value = 6052.9165
$ws.Range("I126").Value = 6065.5
$ws.Range("J126").Value = 5990
$ws.Range("K126").Value = 18196.5
$ws.Range("L126").Value = 17970
$ws.Range("M126").Value = -15726.5
$ws.Range("N126").Value = -22910
$ws.Range("H132").Value = 3451.5
$ws.Range("I132").Value = 2779.9
$ws.Range("K132").Value = 8339.700000000001
$ws.Range("M132").Value = -5809.700000000001
$ws.Range("H136").Value = 36691.93
$ws.Range("I136").Value = 1751.3889
$ws.Range("J136").Value = 93867.37
$ws.Range("K136").Value = 5254.1667
$ws.Range("L136").Value = 281602.11
$ws.Range("M136").Value = -2704.1667
$ws.Range("N136").Value = -286702.11
